$wb = $excel.ActiveWorkbook

# --- Rename the hidden defined name backing the "Range" worksheet
#     connection (بورداصلیD1E31 -> بورداصلیD1E311) ---
foreach ($n in $wb.Names) {
    if ($n.Name -eq "_xlcn.WorksheetConnection_بورداصلیD1E31") {
        $n.Name = "_xlcn.WorksheetConnection_بورداصلیD1E311"
    }
}

# --- Sheet "ماژول ها" (2nd sheet): add a new row with the new task text ---
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("A9").Value = "ماژول مدیریت نود"
$ws2.Rows.Item(9).RowHeight = 26.25

# Move the selection in that sheet down to A10 (below the newly added row)
$ws2.Range("A10").Select()

# Make "ماژول ها" the active/selected tab (was "بورد اصلی")
$ws2.Activate()
